$d = $word.ActiveDocument

# The "<id>...</id>" tag paragraph currently has its contents split over
# several runs:
#   "<id>"  (Courier New, color 7f6000, sz 18/szCs 18)
#   "p14"   (color 000000)
#   "8"     (no explicit color)
#   "v_1"   (color 000000)
#   "</id>" (Courier New, color 7f6000, sz 18/szCs 18)
# Re-join all of that into a single run reading "<id>p148v_1</id>" that
# keeps the opening/closing tag's original formatting, leaving the
# paragraph's other (untouched) runs exactly as they were.

# Locate the whole current (split) text; Word's Find matches straight
# across run boundaries, so this gives us the full span to collapse.
$whole = $d.Content.Duplicate
$whole.Find.Execute("<id>p148v_1</id>", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null

$tagStart = $whole.Start
$tagEnd = $tagStart + 4   # just past the literal "<id>" (4 characters)

# Delete everything from right after the opening tag through the end of the
# closing tag ("p148" + "8" + "v_1" + "</id>" worth of runs). This leaves the
# run holding "<id>" (and whatever follows the old "</id>" run) untouched.
$d.Range($tagEnd, $whole.End).Delete()

# Re-open the (now collapsed) range sitting right at the end of the "<id>"
# run and type the remaining text into it; typing directly after a run with
# matching position simply extends that run (and its formatting) instead of
# starting a differently formatted one.
$insertionPoint = $d.Range($tagEnd, $tagEnd)
$insertionPoint.InsertAfter("p148v_1</id>")
